# Update the "PropertyShapes (properties)" sheet so the datatype / class /
# node columns document that multiple values get auto-wrapped into sh:or,
# and the machine-readable header row carries the separator/wrapper hint.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PropertyShapes (properties)")

# Row 5 = human-readable column descriptions
$ws.Range("J5").Value = "For literal values, the expected datatype of the values. If you put more than one, this will be automatically wrapped into a sh:or."
$ws.Range("K5").Value = "Expected class that the values of the predicate/path must have.  If you put more than one, this will be automatically wrapped into a sh:or."
$ws.Range("L5").Value = "If needed, expected shape that the values of the predicate/path must follow. This must be a reference to a URI of NodeShape from the first sheet.  If you put more than one, this will be automatically wrapped into a sh:or."

# Row 7 = machine-readable SHACL predicate mapping
$ws.Range("J7").Value = "sh:datatype(separator=`",`" wrapper=`"sh:or`")"
$ws.Range("K7").Value = "sh:class(separator=`",`" wrapper=`"sh:or`")"
$ws.Range("L7").Value = "sh:node(separator=`",`" wrapper=`"sh:or`")"

# The longer wrapped descriptions need a taller header row to stay readable.
$ws.Rows.Item(5).RowHeight = 102
